$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume/Number and report week dates) ---
$ws.Range("A8").Value = "Volume 31   Number  40"
$ws.Range("C9").Value = "Report Covering the Week  9/30/2024  Through  10/6/2024"

# --- Cells whose type/style change: copy a same-style template cell first, then set the real value ---
# Row 22
$ws.Range("I14").Copy($ws.Range("C22"))
$ws.Range("C22").Value = 1
$ws.Range("I14").Copy($ws.Range("D22"))
$ws.Range("D22").Value = 1
$ws.Range("L14").Copy($ws.Range("E22"))
$ws.Range("E22").Value = 0
# Row 23
$ws.Range("C14").Copy($ws.Range("D23"))
$ws.Range("E14").Copy($ws.Range("E23"))
# Row 27
$ws.Range("I14").Copy($ws.Range("D27"))
$ws.Range("D27").Value = 1
$ws.Range("L14").Copy($ws.Range("E27"))
$ws.Range("E27").Value = 0
# Row 29
$ws.Range("C14").Copy($ws.Range("C29"))
$ws.Range("C14").Copy($ws.Range("D29"))
$ws.Range("E14").Copy($ws.Range("E29"))
# Row 30
$ws.Range("C14").Copy($ws.Range("C30"))
$ws.Range("C14").Copy($ws.Range("D30"))
$ws.Range("E14").Copy($ws.Range("E30"))
# Row 31
$ws.Range("C14").Copy($ws.Range("C31"))

# --- Plain numeric value updates (style/type unchanged) ---
# Row 15
$ws.Range("C15").Value = 1
$ws.Range("F15").Value = 6
$ws.Range("H15").Value = 20
$ws.Range("I15").Value = 32
$ws.Range("K15").Value = 60
$ws.Range("L15").Value = 10.344827586206
$ws.Range("M15").Value = 100
$ws.Range("N15").Value = -43.859649122807
# Row 16
$ws.Range("C16").Value = 7
$ws.Range("D16").Value = 8
$ws.Range("E16").Value = -12.5
$ws.Range("F16").Value = 33
$ws.Range("G16").Value = 32
$ws.Range("H16").Value = 3.125
$ws.Range("I16").Value = 332
$ws.Range("J16").Value = 284
$ws.Range("K16").Value = 16.901408450704
$ws.Range("L16").Value = 9.570957095709
$ws.Range("M16").Value = -4.871060171919
$ws.Range("N16").Value = -73.692551505546
# Row 17
$ws.Range("C17").Value = 16
$ws.Range("D17").Value = 6
$ws.Range("E17").Value = 166.666666666667
$ws.Range("F17").Value = 60
$ws.Range("G17").Value = 51
$ws.Range("H17").Value = 17.647058823529
$ws.Range("I17").Value = 605
$ws.Range("J17").Value = 504
$ws.Range("K17").Value = 20.039682539682
$ws.Range("L17").Value = 29.273504273504
$ws.Range("M17").Value = 130.038022813688
$ws.Range("N17").Value = -8.472012102874
# Row 18
$ws.Range("C18").Value = 2
$ws.Range("D18").Value = 1
$ws.Range("E18").Value = 100
$ws.Range("F18").Value = 22
$ws.Range("G18").Value = 14
$ws.Range("H18").Value = 57.142857142857
$ws.Range("I18").Value = 179
$ws.Range("J18").Value = 126
$ws.Range("K18").Value = 42.063492063492
$ws.Range("L18").Value = 11.180124223602
$ws.Range("M18").Value = 0
$ws.Range("N18").Value = -79.796839729119
# Row 19
$ws.Range("D19").Value = 10
$ws.Range("E19").Value = 0
$ws.Range("F19").Value = 42
$ws.Range("G19").Value = 54
$ws.Range("H19").Value = -22.222222222222
$ws.Range("I19").Value = 459
$ws.Range("J19").Value = 442
$ws.Range("K19").Value = 3.846153846153
$ws.Range("L19").Value = 1.773835920177
$ws.Range("M19").Value = 51.986754966887
$ws.Range("N19").Value = -50.751072961373
# Row 20
$ws.Range("C20").Value = 7
$ws.Range("D20").Value = 7
$ws.Range("E20").Value = 0
$ws.Range("F20").Value = 14
$ws.Range("G20").Value = 19
$ws.Range("H20").Value = -26.315789473684
$ws.Range("I20").Value = 152
$ws.Range("J20").Value = 181
$ws.Range("K20").Value = -16.022099447513
$ws.Range("L20").Value = -12.138728323699
$ws.Range("M20").Value = 13.432835820895
$ws.Range("N20").Value = -87.151310228233
# Row 21
$ws.Range("C21").Value = 43
$ws.Range("D21").Value = 32
$ws.Range("E21").Value = 34.375
$ws.Range("F21").Value = 177
$ws.Range("G21").Value = 175
$ws.Range("H21").Value = 1.142857142857
$ws.Range("I21").Value = 1764
$ws.Range("J21").Value = 1557
$ws.Range("K21").Value = 13.294797687861
$ws.Range("L21").Value = 10.595611285266
$ws.Range("M21").Value = 40.669856459330
$ws.Range("N21").Value = -64.748201438848
# Row 22
$ws.Range("G22").Value = 3
$ws.Range("H22").Value = 33.333333333333
$ws.Range("I22").Value = 26
$ws.Range("J22").Value = 16
$ws.Range("K22").Value = 62.5
$ws.Range("L22").Value = 52.941176470588
$ws.Range("M22").Value = 36.842105263157
# Row 23
$ws.Range("G23").Value = 4
$ws.Range("H23").Value = -25
$ws.Range("I23").Value = 37
$ws.Range("K23").Value = 2.777777777777
$ws.Range("L23").Value = -15.909090909090
$ws.Range("M23").Value = 19.354838709677
# Row 24
$ws.Range("C24").Value = 63
$ws.Range("D24").Value = 25
$ws.Range("E24").Value = 152
$ws.Range("F24").Value = 228
$ws.Range("G24").Value = 116
$ws.Range("H24").Value = 96.551724137931
$ws.Range("I24").Value = 1595
$ws.Range("J24").Value = 1185
$ws.Range("K24").Value = 34.599156118143
$ws.Range("L24").Value = 36.675235646958
$ws.Range("M24").Value = 80.429864253393
# Row 25
$ws.Range("C25").Value = 38
$ws.Range("D25").Value = 16
$ws.Range("E25").Value = 137.5
$ws.Range("F25").Value = 143
$ws.Range("G25").Value = 61
$ws.Range("H25").Value = 134.426229508197
$ws.Range("I25").Value = 997
$ws.Range("J25").Value = 498
$ws.Range("K25").Value = 100.200803212851
$ws.Range("L25").Value = 79.963898916967
# Row 26
$ws.Range("C26").Value = 20
$ws.Range("D26").Value = 14
$ws.Range("E26").Value = 42.857142857142
$ws.Range("G26").Value = 72
$ws.Range("H26").Value = 43.055555555555
$ws.Range("I26").Value = 843
$ws.Range("J26").Value = 754
$ws.Range("K26").Value = 11.803713527851
$ws.Range("L26").Value = 42.881355932203
$ws.Range("M26").Value = 35.530546623794
# Row 27
$ws.Range("C27").Value = 1
$ws.Range("G27").Value = 7
$ws.Range("H27").Value = 14.285714285714
$ws.Range("I27").Value = 45
$ws.Range("J27").Value = 30
$ws.Range("K27").Value = 50
$ws.Range("L27").Value = 4.651162790697
# Row 28
$ws.Range("D28").Value = 2
$ws.Range("E28").Value = -50
$ws.Range("G28").Value = 9
$ws.Range("H28").Value = -44.444444444444
$ws.Range("I28").Value = 98
$ws.Range("J28").Value = 68
$ws.Range("K28").Value = 44.117647058823
$ws.Range("L28").Value = 55.555555555555
# Row 29
$ws.Range("G29").Value = 2
$ws.Range("H29").Value = 100
$ws.Range("L29").Value = -13.043478260869
$ws.Range("N29").Value = -82.905982905982
# Row 30
$ws.Range("G30").Value = 2
$ws.Range("H30").Value = 50
$ws.Range("L30").Value = -20
$ws.Range("N30").Value = -85.185185185185
# Row 31
$ws.Range("G31").Value = 1
$ws.Range("H31").Value = 0
